$wb = $excel.ActiveWorkbook

# Remove the extra "Sheet1" worksheet, keeping only "Create WO"
$excel.DisplayAlerts = $false
$wb.Worksheets("Sheet1").Delete()

# Update the lot-track item name on the remaining sheet
$ws = $wb.Worksheets("Create WO")
$ws.Range("B2").Value = "Pro-SYDATA1 (Lot track)"

# Move the active selection to B2 to match the saved view state
$ws.Range("B2").Select()
